$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) with the new set of emotion categories.
# Three old categories are dropped ("cười lăn lộn", "cười", "mặt nhăn")
# and three new ones are introduced ("há hốc mồm", "ha ha", "nhăn nhó").
$ws.Range("A1").Value = "lè lưỡi"
$ws.Range("B1").Value = "há hốc mồm"
$ws.Range("C1").Value = "cười mĩm"
$ws.Range("D1").Value = "ha ha"
$ws.Range("E1").Value = "cười ra nước mắt"
$ws.Range("F1").Value = "cảm thấy buồn"
$ws.Range("G1").Value = "muốn khóc"
$ws.Range("H1").Value = "nhăn nhó"
$ws.Range("I1").Value = "rối rắm"
$ws.Range("J1").Value = "cạn lời"
$ws.Range("K1").Value = "cười híp mắt"

# Move the active selection to K1 as in the saved workbook.
$ws.Range("K1").Select()
